# Updates the "dSF" column (F) values for a handful of rows in the
# active worksheet to reflect re-pulled / re-pushed source data and an
# updated mean calculation (per the commit message).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    4  = -2
    5  = 0
    6  = 4
    10 = 1
    11 = 1
    14 = 0
    26 = -1
    27 = 1
    29 = -1
    33 = 2
    39 = 2
    42 = 0
    45 = 0
    57 = 0
    62 = -3
    67 = -2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
